# "moved dates around, pushed midterm to week 6"
#
# Midterm 1 used to be due on the very first due-date in the list (row 9,
# 2024-02-07). It's being pushed back to week 6, landing on the due date
# that the "Programming Project 4" row currently has (row 13, 2024-02-21).
# That means a new row has to be inserted right before "Module 7
# Programming Problems" (row 14) to hold the rescheduled Midterm 1, row 9
# becomes "Programming Project 3" (the assessment that now leads the
# schedule), and every row from the old row 13 onward shifts down by one -
# carrying both their data and their sequential "number" value with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 is no longer Midterm 1 - the next assessment in the schedule
# (Programming Project 3) moves up into its place; its due date is kept.
$ws.Range("B9").Value = "Programming Project 3"

# Insert a fresh row at 13 - this pushes the old row 13 ("Programming
# Project 4" / 2024-02-21) and everything after it down by one.
$ws.Rows("13").Insert()

# The new row 13 is the rescheduled Midterm 1, due on the date that used
# to belong to the row it displaced (now sitting at row 14).
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Midterm 1"
$ws.Range("C13").Value = $ws.Range("C14").Value()

# The "number" column is a plain sequential count, not a formula, so the
# insert didn't renumber it - fix up every row after the inserted one.
for ($r = 14; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Match the author's final cursor position.
$null = $ws.Range("D35").Select()
